$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6162230000000001
$ws.Range("H2").Value = 1.848669
$ws.Range("I2").Value = 0.01637276483811898
$ws.Range("J2").Value = 0.01637276483811898
$ws.Range("M2").Value = 9.101794333333332
$ws.Range("N2").Value = 27.305383
$ws.Range("O2").Value = 0.1526015110517656
$ws.Range("P2").Value = 0.1526015110517656
$ws.Range("Q2").Value = 5.608735009469667
$ws.Range("R2").Value = 50.478615085227
$ws.Range("S2").Value = 0.002498508654392173
$ws.Range("T2").Value = 0.002498508654392173

$ws.Range("G3").Value = 0.6162230000000001
$ws.Range("H3").Value = 1.848669
$ws.Range("I3").Value = 0.01637276483811898
$ws.Range("J3").Value = 0.01637276483811898
$ws.Range("O3").Value = 0.5991759712230392
$ws.Range("P3").Value = 0.5991759712230392
$ws.Range("Q3").Value = 22.02218853187934
$ws.Range("R3").Value = 198.199696786914
$ws.Range("S3").Value = 0.009810167273486367
$ws.Range("T3").Value = 0.009810167273486367

$ws.Range("G4").Value = 0.6162230000000001
$ws.Range("H4").Value = 1.848669
$ws.Range("I4").Value = 0.01637276483811898
$ws.Range("J4").Value = 0.01637276483811898
$ws.Range("O4").Value = 0.2482225177251951
$ws.Range("P4").Value = 0.2482225177251951
$ws.Range("Q4").Value = 9.123201439543669
$ws.Range("R4").Value = 82.10881295589301
$ws.Range("S4").Value = 0.004064088910240441
$ws.Range("T4").Value = 0.00406408891024044

$ws.Range("I5").Value = 0.7731994397120591
$ws.Range("J5").Value = 0.7731994397120592
$ws.Range("M5").Value = 9.101794333333332
$ws.Range("N5").Value = 27.305383
$ws.Range("O5").Value = 0.1526015110517656
$ws.Range("P5").Value = 0.1526015110517656
$ws.Range("Q5").Value = 264.8710104672574
$ws.Range("R5").Value = 2383.839094205317
$ws.Range("S5").Value = 0.1179914028444387
$ws.Range("T5").Value = 0.1179914028444388

$ws.Range("I6").Value = 0.7731994397120591
$ws.Range("J6").Value = 0.7731994397120592
$ws.Range("O6").Value = 0.5991759712230392
$ws.Range("P6").Value = 0.5991759712230392
$ws.Range("S6").Value = 0.4632825252385828
$ws.Range("T6").Value = 0.4632825252385828

$ws.Range("I7").Value = 0.7731994397120591
$ws.Range("J7").Value = 0.7731994397120592
$ws.Range("O7").Value = 0.2482225177251951
$ws.Range("P7").Value = 0.2482225177251951
$ws.Range("S7").Value = 0.1919255116290375
$ws.Range("T7").Value = 0.1919255116290375

$ws.Range("G8").Value = 7.919886999999999
$ws.Range("I8").Value = 0.2104277954498219
$ws.Range("J8").Value = 0.2104277954498219
$ws.Range("M8").Value = 9.101794333333332
$ws.Range("N8").Value = 27.305383
$ws.Range("O8").Value = 0.1526015110517656
$ws.Range("P8").Value = 0.1526015110517656
$ws.Range("Q8").Value = 72.08518261724032
$ws.Range("R8").Value = 648.7666435551629
$ws.Range("S8").Value = 0.03211159955293467
$ws.Range("T8").Value = 0.03211159955293467

$ws.Range("G9").Value = 7.919886999999999
$ws.Range("I9").Value = 0.2104277954498219
$ws.Range("J9").Value = 0.2104277954498219
$ws.Range("O9").Value = 0.5991759712230392
$ws.Range("P9").Value = 0.5991759712230392
$ws.Range("S9").Value = 0.1260832787109701
$ws.Range("T9").Value = 0.1260832787109701

$ws.Range("G10").Value = 7.919886999999999
$ws.Range("I10").Value = 0.2104277954498219
$ws.Range("J10").Value = 0.2104277954498219
$ws.Range("O10").Value = 0.2482225177251951
$ws.Range("P10").Value = 0.2482225177251951
$ws.Range("Q10").Value = 117.2541831113463
$ws.Range("S10").Value = 0.05223291718591715
$ws.Range("T10").Value = 0.05223291718591715

